# Horarios actualizados Linea 141 - 252
# Updates the scraped-schedule snapshot (new scrape at 04:40:33) across the
# three sheets of the workbook: existing rows shift their values down as new,
# earlier-arriving buses are inserted at the top of each sheet's data block,
# and newly-scraped rows are appended at the bottom. Header rows (1/3/5) and
# totals are refreshed to match.
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = 'Última actualización: 04:40:33'
$ws.Range("A3").Value = 'Total filas: 19'
$ws.Range("A7").Value = '04:40:33'
$ws.Range("B7").Value = '04:46'
$ws.Range("C7").Value = '215A_EL PATO'
$ws.Range("D7").Value = 6
$ws.Range("B8").Value = '04:53'
$ws.Range("C8").Value = '11_ETCHEVERRY'
$ws.Range("D8").Value = 35
$ws.Range("B9").Value = '05:16'
$ws.Range("C9").Value = '17_ROMERO'
$ws.Range("D9").Value = 58
$ws.Range("A10").Value = '04:40:33'
$ws.Range("B10").Value = '05:16'
$ws.Range("C10").Value = '14_ABASTO'
$ws.Range("D10").Value = 36
$ws.Range("B11").Value = '05:21'
$ws.Range("C11").Value = '23_HERNANDEZ'
$ws.Range("D11").Value = 63
$ws.Range("A12").Value = '04:40:33'
$ws.Range("B12").Value = '05:22'
$ws.Range("C12").Value = '23_HERNANDEZ'
$ws.Range("D12").Value = 42
$ws.Range("B13").Value = '05:34'
$ws.Range("C13").Value = '215B_EL PATO'
$ws.Range("D13").Value = 76
$ws.Range("B14").Value = '05:37'
$ws.Range("C14").Value = '14_ABASTO'
$ws.Range("D14").Value = 79
$ws.Range("B15").Value = '05:46'
$ws.Range("C15").Value = '15_ABASTO'
$ws.Range("D15").Value = 88
$ws.Range("A16").Value = '04:40:33'
$ws.Range("B16").Value = '06:04'
$ws.Range("C16").Value = '16_SANTA ANA'
$ws.Range("D16").Value = 84
$ws.Range("E16").Value = 'LP1912'
$ws.Range("A17").Value = '04:18:53'
$ws.Range("B17").Value = '06:07'
$ws.Range("C17").Value = '16_SANTA ANA'
$ws.Range("D17").Value = 109
$ws.Range("E17").Value = 'LP1912'
$ws.Range("A18").Value = '04:18:53'
$ws.Range("B18").Value = '06:11'
$ws.Range("C18").Value = '215A_EL PATO'
$ws.Range("D18").Value = 113
$ws.Range("E18").Value = 'LP1912'
$ws.Range("A19").Value = '04:18:53'
$ws.Range("B19").Value = '06:13'
$ws.Range("C19").Value = '225_HARAS DEL SUR'
$ws.Range("D19").Value = 115
$ws.Range("E19").Value = 'LP1912'
$ws.Range("A20").Value = '04:40:33'
$ws.Range("B20").Value = '06:14'
$ws.Range("C20").Value = '225_HARAS DEL SUR'
$ws.Range("D20").Value = 94
$ws.Range("E20").Value = 'LP1912'
$ws.Range("A21").Value = '04:40:33'
$ws.Range("B21").Value = '06:21'
$ws.Range("C21").Value = '26_HERNANDEZ'
$ws.Range("D21").Value = 101
$ws.Range("E21").Value = 'LP1912'
$ws.Range("A22").Value = '04:40:33'
$ws.Range("B22").Value = '06:27'
$ws.Range("C22").Value = '23_HERNANDEZ'
$ws.Range("D22").Value = 107
$ws.Range("E22").Value = 'LP1912'
$ws.Range("A23").Value = '04:40:33'
$ws.Range("B23").Value = '06:29'
$ws.Range("C23").Value = '86_EST CHICA-ESC AGRARIA'
$ws.Range("D23").Value = 109
$ws.Range("E23").Value = 'LP1912'
$ws.Range("A24").Value = '04:40:33'
$ws.Range("B24").Value = '06:31'
$ws.Range("C24").Value = '16_SANTA ANA'
$ws.Range("D24").Value = 111
$ws.Range("E24").Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = 'Última actualización: 04:40:33'
$ws.Range("A3").Value = 'Total filas: 4'
$ws.Range("A7").Value = '04:40:33'
$ws.Range("B7").Value = '04:46'
$ws.Range("C7").Value = '215A_EL PATO'
$ws.Range("D7").Value = 6
$ws.Range("B8").Value = '05:34'
$ws.Range("C8").Value = '215B_EL PATO'
$ws.Range("D8").Value = 76
$ws.Range("A9").Value = '04:18:53'
$ws.Range("B9").Value = '06:11'
$ws.Range("C9").Value = '215A_EL PATO'
$ws.Range("D9").Value = 113
$ws.Range("E9").Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = 'Última actualización: 04:40:33'
$ws.Range("A3").Value = 'Total filas: 5'
$ws.Range("A7").Value = '04:40:33'
$ws.Range("B7").Value = '05:44'
$ws.Range("D7").Value = 64
$ws.Range("A8").Value = '04:18:53'
$ws.Range("B8").Value = '06:08'
$ws.Range("C8").Value = '215A_LA PLATA'
$ws.Range("D8").Value = 110
$ws.Range("E8").Value = 'L6173'
$ws.Range("A9").Value = '04:40:33'
$ws.Range("B9").Value = '06:09'
$ws.Range("C9").Value = '215A_LA PLATA'
$ws.Range("D9").Value = 89
$ws.Range("E9").Value = 'L6173'
$ws.Range("A10").Value = '04:40:33'
$ws.Range("B10").Value = '06:33'
$ws.Range("C10").Value = '215C_LA PLATA'
$ws.Range("D10").Value = 113
$ws.Range("E10").Value = 'L6203'
